$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) and the "SC 92" row (originally row 28).
# Deleting row 26 first shifts "SC 92" up to row 27, so delete that row next.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Fix individual cell values that changed in the remaining data.
$ws.Range("F19").Value = 17.81
$ws.Range("F21").ClearContents()
$ws.Range("F23").Value = 16.48
$ws.Range("F27").ClearContents()
$ws.Range("F33").Value = 17.53

$wb.Save()
